# Regenerate orders with updated distance/sizes.
# The stimulus-order workbook encodes trial "Distance" (D..) and "Size" (S..)
# codes inside several text columns (Condition, Filename_Left, Filename_Right,
# Distance, Size). This commit renumbers those codes:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# Every occurrence of these codes - whether the whole cell value (e.g. the
# "Distance"/"Size" lookup columns) or embedded inside a longer token (e.g.
# "Face17_D80_S25" or "Face17_D80_S25_l.png") - must be updated consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @(
    @("D64", "D69"),
    @("D80", "D86"),
    @("D51", "D55"),
    @("S30", "S31")
)

$rng = $ws.UsedRange

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replaceWith = $pair[1]
    # LookAt:=xlPart (2) so substrings embedded within larger tokens
    # (e.g. Face17_D80_S25_l.png) are matched too; MatchCase:=True so we
    # only touch the exact uppercase codes.
    $rng.Replace($find, $replaceWith, 2, 1, $false, $false, $true, $false)
}
